# Apply a cyclic rotation of species-record data across rows 26, 27 and 28:
#   old row 26 data -> row 28
#   old row 27 data -> row 26
#   old row 28 data -> row 27
# Only columns A, B, E, F, G, H, Q, R change; all other columns stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 (was row 27's content)
$ws.Range("A26").Value = 111541116
$ws.Range("B26").Value = 108219
$ws.Range("E26").Value = 219711
$ws.Range("F26").Value = "Sårläka"
$ws.Range("G26").Value = "Sanicula europaea"
$ws.Range("H26").Value = "L."
$ws.Range("Q26").Value = 693830.7552326696
$ws.Range("R26").Value = 6552178.401404973

# Row 27 (was row 28's content)
$ws.Range("A27").Value = 111541117
$ws.Range("B27").Value = 108219
$ws.Range("E27").Value = 219711
$ws.Range("F27").Value = "Sårläka"
$ws.Range("G27").Value = "Sanicula europaea"
$ws.Range("H27").Value = "L."
$ws.Range("Q27").Value = 693809.5100469354
$ws.Range("R27").Value = 6552200.504896822

# Row 28 (was row 26's content)
$ws.Range("A28").Value = 111541130
$ws.Range("B28").Value = 98535
$ws.Range("E28").Value = 222498
$ws.Range("F28").Value = "Blåsippa"
$ws.Range("G28").Value = "Hepatica nobilis"
$ws.Range("H28").Value = "Schreb."
$ws.Range("Q28").Value = 693830.8333423812
$ws.Range("R28").Value = 6552176.860022029
